# Case_1_76 (380 kV) res_line/loading_percent.xlsx results refresh.
# Sets the updated loading-percent values (rows 2-25 = time steps 0-23,
# columns B,D,E,F,G,H,I,J,K,M,O = line IDs); columns A,C,L,N are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 6.072429246268456
$ws.Cells.Item(2, 4).Value = 10.01628844949731
$ws.Cells.Item(2, 5).Value = 14.25931973078596
$ws.Cells.Item(2, 6).Value = 29.487398179572
$ws.Cells.Item(2, 7).Value = 28.53853682931905
$ws.Cells.Item(2, 8).Value = 14.32415445241967
$ws.Cells.Item(2, 9).Value = 20.3434681855076
$ws.Cells.Item(2, 10).Value = 10.1455275034992
$ws.Cells.Item(2, 11).Value = 9.96213452487169
$ws.Cells.Item(2, 13).Value = 14.81273801262293
$ws.Cells.Item(2, 15).Value = 21.75769422644195
$ws.Cells.Item(3, 2).Value = 5.919042429191708
$ws.Cells.Item(3, 4).Value = 9.992922786501298
$ws.Cells.Item(3, 5).Value = 14.27345791830484
$ws.Cells.Item(3, 6).Value = 29.56934621847424
$ws.Cells.Item(3, 7).Value = 28.64038434121811
$ws.Cells.Item(3, 8).Value = 14.37551779368246
$ws.Cells.Item(3, 9).Value = 20.46028066720997
$ws.Cells.Item(3, 10).Value = 10.1725361870417
$ws.Cells.Item(3, 11).Value = 9.557865457257082
$ws.Cells.Item(3, 13).Value = 14.65813713836883
$ws.Cells.Item(3, 15).Value = 21.84430683723154
$ws.Cells.Item(4, 2).Value = 5.823446154114829
$ws.Cells.Item(4, 4).Value = 9.980097968653466
$ws.Cells.Item(4, 5).Value = 14.28463725997879
$ws.Cells.Item(4, 6).Value = 29.62683088936815
$ws.Cells.Item(4, 7).Value = 28.71238466286596
$ws.Cells.Item(4, 8).Value = 14.40935069746151
$ws.Cells.Item(4, 9).Value = 20.5357378155487
$ws.Cells.Item(4, 10).Value = 10.19033606055065
$ws.Cells.Item(4, 11).Value = 9.300647029068831
$ws.Cells.Item(4, 13).Value = 14.5639179455439
$ws.Cells.Item(4, 15).Value = 21.90221542536285
$ws.Cells.Item(5, 2).Value = 5.784189063495451
$ws.Cells.Item(5, 4).Value = 9.97525829078082
$ws.Cells.Item(5, 5).Value = 14.28982154933821
$ws.Cells.Item(5, 6).Value = 29.65205483312161
$ws.Cells.Item(5, 7).Value = 28.74409466151548
$ws.Cells.Item(5, 8).Value = 14.42371527647126
$ws.Cells.Item(5, 9).Value = 20.56742866593846
$ws.Cells.Item(5, 10).Value = 10.19789588641043
$ws.Cells.Item(5, 11).Value = 9.193699190070346
$ws.Cells.Item(5, 13).Value = 14.52573517426814
$ws.Cells.Item(5, 15).Value = 21.92700063488296
$ws.Cells.Item(6, 2).Value = 5.777653955650037
$ws.Cells.Item(6, 4).Value = 9.974478119267745
$ws.Cells.Item(6, 5).Value = 14.29072037091546
$ws.Cells.Item(6, 6).Value = 29.65635174846051
$ws.Cells.Item(6, 7).Value = 28.749502850327
$ws.Cells.Item(6, 8).Value = 14.42613538169875
$ws.Cells.Item(6, 9).Value = 20.57274784399664
$ws.Cells.Item(6, 10).Value = 10.19916969580112
$ws.Cells.Item(6, 11).Value = 9.175816193634441
$ws.Cells.Item(6, 13).Value = 14.5194087841298
$ws.Cells.Item(6, 15).Value = 21.93118783443024
$ws.Cells.Item(7, 2).Value = 5.822917862161446
$ws.Cells.Item(7, 4).Value = 9.980031129082027
$ws.Cells.Item(7, 5).Value = 14.28470463150907
$ws.Cells.Item(7, 6).Value = 29.62716379193532
$ws.Cells.Item(7, 7).Value = 28.71280273675121
$ws.Cells.Item(7, 8).Value = 14.40954208514283
$ws.Cells.Item(7, 9).Value = 20.53616139414114
$ws.Cells.Item(7, 10).Value = 10.19043677461169
$ws.Cells.Item(7, 11).Value = 9.299213126532436
$ws.Cells.Item(7, 13).Value = 14.56340209412468
$ws.Cells.Item(7, 15).Value = 21.90254488386282
$ws.Cells.Item(8, 2).Value = 6.019867637907661
$ws.Cells.Item(8, 4).Value = 10.00791853038011
$ws.Cells.Item(8, 5).Value = 14.26367629279416
$ws.Cells.Item(8, 6).Value = 29.51416376520935
$ws.Cells.Item(8, 7).Value = 28.57168291665021
$ws.Cells.Item(8, 8).Value = 14.34138816862283
$ws.Cells.Item(8, 9).Value = 20.38297171575575
$ws.Cells.Item(8, 10).Value = 10.15458783842349
$ws.Cells.Item(8, 11).Value = 9.82468412601267
$ws.Cells.Item(8, 13).Value = 14.75930593104336
$ws.Cells.Item(8, 15).Value = 21.78657578997064
$ws.Cells.Item(9, 2).Value = 6.392575369561015
$ws.Cells.Item(9, 4).Value = 10.07449823020375
$ws.Cells.Item(9, 5).Value = 14.24224226094719
$ws.Cells.Item(9, 6).Value = 29.34961517906467
$ws.Cells.Item(9, 7).Value = 28.3705101349591
$ws.Cells.Item(9, 8).Value = 14.22594640950124
$ws.Cells.Item(9, 9).Value = 20.11207626372855
$ws.Cells.Item(9, 10).Value = 10.09392427908987
$ws.Cells.Item(9, 11).Value = 10.77873812865385
$ws.Cells.Item(9, 13).Value = 15.14754317718768
$ws.Cells.Item(9, 15).Value = 21.59676105722122
$ws.Cells.Item(10, 2).Value = 6.655353907525352
$ws.Cells.Item(10, 4).Value = 10.13040312862296
$ws.Cells.Item(10, 5).Value = 14.23852962989466
$ws.Cells.Item(10, 6).Value = 29.26370550165017
$ws.Cells.Item(10, 7).Value = 28.26935409109235
$ws.Cells.Item(10, 8).Value = 14.15222194245262
$ws.Cells.Item(10, 9).Value = 19.930874132787
$ws.Cells.Item(10, 10).Value = 10.05520751419566
$ws.Cells.Item(10, 11).Value = 11.4803650648151
$ws.Cells.Item(10, 13).Value = 15.43320552644317
$ws.Cells.Item(10, 15).Value = 21.48034072730661
$ws.Cells.Item(11, 2).Value = 6.771972678508733
$ws.Cells.Item(11, 4).Value = 10.15729182446812
$ws.Cells.Item(11, 5).Value = 14.23944233970081
$ws.Cells.Item(11, 6).Value = 29.23225606129325
$ws.Cells.Item(11, 7).Value = 28.23357215662212
$ws.Cells.Item(11, 8).Value = 14.12108922026959
$ws.Cells.Item(11, 9).Value = 19.85227614967103
$ws.Cells.Item(11, 10).Value = 10.03886057705459
$ws.Cells.Item(11, 11).Value = 11.78743704487558
$ws.Cells.Item(11, 13).Value = 15.5628194336827
$ws.Cells.Item(11, 15).Value = 21.4324041845499
$ws.Cells.Item(12, 2).Value = 6.815673995721276
$ws.Cells.Item(12, 4).Value = 10.16767778832577
$ws.Cells.Item(12, 5).Value = 14.24016073221159
$ws.Cells.Item(12, 6).Value = 29.22144662795383
$ws.Cells.Item(12, 7).Value = 28.2215016485077
$ws.Cells.Item(12, 8).Value = 14.10964582491724
$ws.Cells.Item(12, 9).Value = 19.82306155387946
$ws.Cells.Item(12, 10).Value = 10.03285205983359
$ws.Cells.Item(12, 11).Value = 11.90141503503798
$ws.Cells.Item(12, 13).Value = 15.61181682856017
$ws.Cells.Item(12, 15).Value = 21.41497622291207
$ws.Cells.Item(13, 2).Value = 6.806283175507284
$ws.Cells.Item(13, 4).Value = 10.16543201551013
$ws.Cells.Item(13, 5).Value = 14.2399894579551
$ws.Cells.Item(13, 6).Value = 29.22372568598301
$ws.Cells.Item(13, 7).Value = 28.22403534315409
$ws.Cells.Item(13, 8).Value = 14.11209498008489
$ws.Cells.Item(13, 9).Value = 19.82932906931752
$ws.Cells.Item(13, 10).Value = 10.03413802399592
$ws.Cells.Item(13, 11).Value = 11.87697085212268
$ws.Cells.Item(13, 13).Value = 15.60126869634748
$ws.Cells.Item(13, 15).Value = 21.41869738692563
$ws.Cells.Item(14, 2).Value = 6.775577440157651
$ws.Cells.Item(14, 4).Value = 10.1581422334037
$ws.Cells.Item(14, 5).Value = 14.23949397980944
$ws.Cells.Item(14, 6).Value = 29.23134470914716
$ws.Cells.Item(14, 7).Value = 28.2325494202366
$ws.Cells.Item(14, 8).Value = 14.12014083265625
$ws.Cells.Item(14, 9).Value = 19.84986166221831
$ws.Cells.Item(14, 10).Value = 10.03836261244059
$ws.Cells.Item(14, 11).Value = 11.79686045448856
$ws.Cells.Item(14, 13).Value = 15.56685237982344
$ws.Cells.Item(14, 15).Value = 21.43095584199675
$ws.Cells.Item(15, 2).Value = 6.756708306562373
$ws.Cells.Item(15, 4).Value = 10.15370339680895
$ws.Cells.Item(15, 5).Value = 14.2392389875663
$ws.Cells.Item(15, 6).Value = 29.23615486936363
$ws.Cells.Item(15, 7).Value = 28.23795740167155
$ws.Cells.Item(15, 8).Value = 14.12511419705352
$ws.Cells.Item(15, 9).Value = 19.86250986599934
$ws.Cells.Item(15, 10).Value = 10.04097395140062
$ws.Cells.Item(15, 11).Value = 11.74748946113062
$ws.Cells.Item(15, 13).Value = 15.54575930973242
$ws.Cells.Item(15, 15).Value = 21.43855892327492
$ws.Cells.Item(16, 2).Value = 6.647670354994008
$ws.Cells.Item(16, 4).Value = 10.12867472696828
$ws.Cells.Item(16, 5).Value = 14.23852222321951
$ws.Cells.Item(16, 6).Value = 29.26591457483415
$ws.Cells.Item(16, 7).Value = 28.27189914179743
$ws.Cells.Item(16, 8).Value = 14.15430493756376
$ws.Cells.Item(16, 9).Value = 19.93608760795266
$ws.Cells.Item(16, 10).Value = 10.05630126753177
$ws.Cells.Item(16, 11).Value = 11.45997688493003
$ws.Cells.Item(16, 13).Value = 15.42472508244434
$ws.Cells.Item(16, 15).Value = 21.48357476857981
$ws.Cells.Item(17, 2).Value = 6.580001739342138
$ws.Cells.Item(17, 4).Value = 10.11368967769508
$ws.Cells.Item(17, 5).Value = 14.23874795553785
$ws.Cells.Item(17, 6).Value = 29.28612761400597
$ws.Cells.Item(17, 7).Value = 28.29534831489391
$ws.Cells.Item(17, 8).Value = 14.17282851142905
$ws.Cells.Item(17, 9).Value = 19.98220496900113
$ws.Cells.Item(17, 10).Value = 10.06602801761094
$ws.Cells.Item(17, 11).Value = 11.27953541365289
$ws.Cells.Item(17, 13).Value = 15.35036255165039
$ws.Cells.Item(17, 15).Value = 21.51247876963882
$ws.Cells.Item(18, 2).Value = 6.540807802897749
$ws.Cells.Item(18, 4).Value = 10.10520820054699
$ws.Cells.Item(18, 5).Value = 14.23912275875843
$ws.Cells.Item(18, 6).Value = 29.29847179811204
$ws.Cells.Item(18, 7).Value = 28.3097985759883
$ws.Cells.Item(18, 8).Value = 14.18370913779388
$ws.Cells.Item(18, 9).Value = 20.00909125284102
$ws.Cells.Item(18, 10).Value = 10.07174171884566
$ws.Cells.Item(18, 11).Value = 11.17427100220154
$ws.Cells.Item(18, 13).Value = 15.30756176765103
$ws.Cells.Item(18, 15).Value = 21.52957630926285
$ws.Cells.Item(19, 2).Value = 6.527491832963282
$ws.Cells.Item(19, 4).Value = 10.10236030970893
$ws.Cells.Item(19, 5).Value = 14.23929177169641
$ws.Cells.Item(19, 6).Value = 29.30277459351053
$ws.Cells.Item(19, 7).Value = 28.31485632025371
$ws.Cells.Item(19, 8).Value = 14.18743200597704
$ws.Cells.Item(19, 9).Value = 20.01825652293025
$ws.Cells.Item(19, 10).Value = 10.07369674907533
$ws.Cells.Item(19, 11).Value = 11.13837758673495
$ws.Cells.Item(19, 13).Value = 15.29306622234862
$ws.Cells.Item(19, 15).Value = 21.53544636601245
$ws.Cells.Item(20, 2).Value = 6.587233705007291
$ws.Cells.Item(20, 4).Value = 10.11527066915948
$ws.Cells.Item(20, 5).Value = 14.23869858175201
$ws.Cells.Item(20, 6).Value = 29.28390155372676
$ws.Cells.Item(20, 7).Value = 28.29275239383712
$ws.Cells.Item(20, 8).Value = 14.17083321784015
$ws.Cells.Item(20, 9).Value = 19.977258372389
$ws.Cells.Item(20, 10).Value = 10.06498026069765
$ws.Cells.Item(20, 11).Value = 11.29889713224701
$ws.Cells.Item(20, 13).Value = 15.35828188059932
$ws.Cells.Item(20, 15).Value = 21.5093529484591
$ws.Cells.Item(21, 2).Value = 6.784609234058026
$ws.Cells.Item(21, 4).Value = 10.16027793520381
$ws.Cells.Item(21, 5).Value = 14.23962940855773
$ws.Cells.Item(21, 6).Value = 29.22907695249887
$ws.Cells.Item(21, 7).Value = 28.23000842456798
$ws.Cells.Item(21, 8).Value = 14.11776818275495
$ws.Cells.Item(21, 9).Value = 19.84381586671344
$ws.Cells.Item(21, 10).Value = 10.03711681887589
$ws.Cells.Item(21, 11).Value = 11.82045362328865
$ws.Cells.Item(21, 13).Value = 15.5769638681513
$ws.Cells.Item(21, 15).Value = 21.42733555499759
$ws.Cells.Item(22, 2).Value = 6.910910158041014
$ws.Cells.Item(22, 4).Value = 10.1908778344857
$ws.Cells.Item(22, 5).Value = 14.24240977984338
$ws.Cells.Item(22, 6).Value = 29.19965709763975
$ws.Cells.Item(22, 7).Value = 28.19762717016007
$ws.Cells.Item(22, 8).Value = 14.08510326943584
$ws.Cells.Item(22, 9).Value = 19.75980101329496
$ws.Cells.Item(22, 10).Value = 10.01996549026448
$ws.Cells.Item(22, 11).Value = 12.14788060955388
$ws.Cells.Item(22, 13).Value = 15.71937583249779
$ws.Cells.Item(22, 15).Value = 21.37795674805955
$ws.Cells.Item(23, 2).Value = 6.843759846374682
$ws.Cells.Item(23, 4).Value = 10.17443962582425
$ws.Cells.Item(23, 5).Value = 14.24072759679141
$ws.Cells.Item(23, 6).Value = 29.21477172601213
$ws.Cells.Item(23, 7).Value = 28.21411814020437
$ws.Cells.Item(23, 8).Value = 14.10235264529824
$ws.Cells.Item(23, 9).Value = 19.80434948274181
$ws.Cells.Item(23, 10).Value = 10.02902266203964
$ws.Cells.Item(23, 11).Value = 11.97436798649763
$ws.Cells.Item(23, 13).Value = 15.64342632613692
$ws.Cells.Item(23, 15).Value = 21.40392391834886
$ws.Cells.Item(24, 2).Value = 6.583965034654788
$ws.Cells.Item(24, 4).Value = 10.11455548626003
$ws.Cells.Item(24, 5).Value = 14.23872014031706
$ws.Cells.Item(24, 6).Value = 29.28490570231326
$ws.Cells.Item(24, 7).Value = 28.29392299191554
$ws.Cells.Item(24, 8).Value = 14.17173456997056
$ws.Cells.Item(24, 9).Value = 19.97949356731265
$ws.Cells.Item(24, 10).Value = 10.06545357260872
$ws.Cells.Item(24, 11).Value = 11.29014844993194
$ws.Cells.Item(24, 13).Value = 15.35470170112316
$ws.Cells.Item(24, 15).Value = 21.51076463644544
$ws.Cells.Item(25, 2).Value = 6.293498163438768
$ws.Cells.Item(25, 4).Value = 10.0552392945411
$ws.Cells.Item(25, 5).Value = 14.24592384285552
$ws.Cells.Item(25, 6).Value = 29.387999865677
$ws.Cells.Item(25, 7).Value = 28.41678019502652
$ws.Cells.Item(25, 8).Value = 14.25522819967775
$ws.Cells.Item(25, 9).Value = 20.18221864859597
$ws.Cells.Item(25, 10).Value = 10.10930602082432
$ws.Cells.Item(25, 11).Value = 10.52945949032062
$ws.Cells.Item(25, 13).Value = 15.04229171714284
$ws.Cells.Item(25, 15).Value = 21.64407365067959
